$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are kept as text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.355.48"
$ws.Range("E2").Value = "  +3.49%  "
$ws.Range("D3").Value = "1.723.81"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "241.95"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4739"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "0.06205"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "1.720.88"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").Value = "0.07081"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("D13").Value = "0.5966"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "4.434"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "76.47"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "26.371.03"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.000006830"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "11.56"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "1.940.11"
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").Value = "4.522"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").Value = "8.754"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "135.34"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").Value = "15.24"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").Value = "1.777"
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "106.95"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").Value = "3.954"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").Value = "3.690"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "0.07810"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "0.04513"
$ws.Range("E33").Value = "  +6.95%  "
$ws.Range("D34").Value = "2.615"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "0.9819"
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D36").Value = "0.6221"
$ws.Range("E36").Value = "  +2.14%  "
$ws.Range("D37").Value = "0.9464"
$ws.Range("E37").Value = "  +10.66%  "
$ws.Range("D38").Value = "114.41"
$ws.Range("E38").Value = "  +19.01%  "
$ws.Range("D39").Value = "2.461"
$ws.Range("E39").Value = "  -5.27%  "
$ws.Range("D40").Value = "1.930"
$ws.Range("E40").Value = "  +4.33%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "5.719"
$ws.Range("E42").Value = "  +18.46%  "
$ws.Range("D43").Value = "0.01486"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "0.3839"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("E45").Value = "  +6.32%  "
$ws.Range("D46").Value = "6.381"
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("D47").Value = "0.05272"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "7.901"
$ws.Range("E48").Value = "  +7.12%  "
$ws.Range("D49").Value = "30.41"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "0.3388"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "1.217"
$ws.Range("E51").Value = "  +2.06%  "
